# Fill in the two evaluator score rows (row 5 and row 6) that were left
# blank in the template. Columns C:AD hold the rubric scores; F/K/O/T/X/AC
# are the per-section subtotal columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row5 = @{
    "C5" = 5;  "D5" = 4;  "E5" = 5;  "F5" = 14;
    "G5" = 5;  "H5" = 5;  "I5" = 4;  "J5" = 4;  "K5" = 18;
    "L5" = 5;  "M5" = 5;  "N5" = 5;  "O5" = 15;
    "P5" = 4;  "Q5" = 5;  "R5" = 5;  "S5" = 5;  "T5" = 19;
    "U5" = 5;  "V5" = 5;  "W5" = 5;  "X5" = 15;
    "Y5" = 5;  "Z5" = 5;  "AA5" = 5; "AB5" = 5; "AC5" = 20;
    "AD5" = 5
}

$row6 = @{
    "C6" = 5;  "D6" = 4;  "E6" = 5;  "F6" = 14;
    "G6" = 4;  "H6" = 4;  "I6" = 5;  "J6" = 5;  "K6" = 18;
    "L6" = 5;  "M6" = 5;  "N6" = 4;  "O6" = 14;
    "P6" = 5;  "Q6" = 5;  "R6" = 4;  "S6" = 4;  "T6" = 18;
    "U6" = 5;  "V6" = 5;  "W6" = 4;  "X6" = 14;
    "Y6" = 5;  "Z6" = 5;  "AA6" = 5; "AB6" = 4; "AC6" = 19;
    "AD6" = 5
}

foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

foreach ($addr in $row6.Keys) {
    $ws.Range($addr).Value = $row6[$addr]
}
